# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet "2022-Q4" right after "总计", holding the
#    quarter's per-fund holdings (cloned from the "2022-Q1" sheet so it
#    inherits identical header/row styling), pushing the older quarter
#    sheets ("2022-Q1","2021-Q4","2021-Q1","2020-Q4") one slot to the right.
# 2) Update the "总计" (totals) sheet: add a new top data row for 2022-Q4
#    and shift the previous rows down, renumbering the index column.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)
$q1_2022 = $wb.Worksheets.Item(2)

# --- Create the new "2022-Q4" sheet by cloning "2022-Q1" (same layout/style) ---
$q1_2022.Copy($null, $total)
$q4_2022 = $wb.Worksheets.Item(2)
$q4_2022.Name = "2022-Q4"

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Row 2: 010116 / 民生加银新兴产业混合A
$q4_2022.Range("A2").Value = 0
Set-TextValue $q4_2022 "B2" "010116"
Set-TextValue $q4_2022 "C2" "民生加银新兴产业混合A"
Set-TextValue $q4_2022 "D2" "6.99"
Set-TextValue $q4_2022 "E2" "85.62"
Set-TextValue $q4_2022 "F2" "4.07"
Set-TextValue $q4_2022 "G2" "0.2845"
$q4_2022.Range("H2").Value = 7

# Row 3: 010117 / 民生加银新兴产业混合C
$q4_2022.Range("A3").Value = 1
Set-TextValue $q4_2022 "B3" "010117"
Set-TextValue $q4_2022 "C3" "民生加银新兴产业混合C"
Set-TextValue $q4_2022 "D3" "0.76"
Set-TextValue $q4_2022 "E3" "85.62"
Set-TextValue $q4_2022 "F3" "4.07"
Set-TextValue $q4_2022 "G3" "0.0309"
$q4_2022.Range("H3").Value = 7

# Row 4: 001252 / 中海进取收益灵活配置混合
$q4_2022.Range("A4").Value = 2
Set-TextValue $q4_2022 "B4" "001252"
Set-TextValue $q4_2022 "C4" "中海进取收益灵活配置混合"
Set-TextValue $q4_2022 "D4" "0.51"
Set-TextValue $q4_2022 "E4" "36.60"
Set-TextValue $q4_2022 "F4" "1.55"
Set-TextValue $q4_2022 "G4" "0.0079"
$q4_2022.Range("H4").Value = 5

# --- Update the "总计" summary sheet: insert a new 2022-Q4 row on top ---
# Shift existing rows down one position (bottom-up so nothing is clobbered).
$total.Range("A6").Value = 4
$total.Range("B6").Value = "2020-Q4"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.14

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q1"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.11

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q4"
$total.Range("C4").Value = 7
$total.Range("D4").Value = 0.08

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q1"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.02

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.32

# Row 6 is brand new - it needs the index column's style (A2:A5 all carry it
# already from the pre-existing rows); clone it onto A6 via a format-only paste.
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)

# Restore the originally-active sheet ("2020-Q4", the last tab) as the
# selected tab, since creating/renaming the new sheet shifts focus to it.
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
